$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.476000000000001
$ws.Range("B9").Value = 6.485000000000001
$ws.Range("D11").Value = -8.316999999999998
$ws.Range("B18").Value = 6.351
$ws.Range("B20").Value = 6.667999999999999
$ws.Range("E21").Value = 13.123
